$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 is "b.md.md" - status columns B3/C3 change
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is "b.md.md" handoff info
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-02-17 02:55:35"

# de-de sheet: row 3 is "b.md.md" handoff info
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$wsDe.Range("D3").Value = "2016-02-17 02:55:45"
